# Update "想去人数" (want-to-go count) figures in column F for the
# 展览 (sheet 1) and 全部类型 (sheet 4) worksheets, matching the refreshed
# scrape output (gh-pages data generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetExhibitions = $wb.Worksheets.Item("展览")
$sheetAllTypes    = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (exhibitions) — row -> new F value
$exhibitionUpdates = @{
    "F3"  = 337
    "F4"  = 193
    "F5"  = 1216
    "F6"  = 441
    "F8"  = 162
    "F10" = 776
    "F24" = 2609
    "F25" = 1374
    "F27" = 18
    "F28" = 327
    "F29" = 394
    "F31" = 797
    "F32" = 1238
    "F36" = 535
    "F37" = 634
    "F38" = 808
}

foreach ($addr in $exhibitionUpdates.Keys) {
    $sheetExhibitions.Range($addr).Value = $exhibitionUpdates[$addr]
}

# Sheet "全部类型" (all types) — row -> new F value
$allTypesUpdates = @{
    "F7"  = 337
    "F8"  = 193
    "F11" = 1216
    "F12" = 441
    "F14" = 162
    "F27" = 2609
    "F29" = 1374
    "F34" = 327
    "F35" = 394
    "F39" = 797
    "F40" = 1238
    "F42" = 535
    "F43" = 634
    "F44" = 808
}

foreach ($addr in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range($addr).Value = $allTypesUpdates[$addr]
}

Write-Output "Updated $($exhibitionUpdates.Count) cells on 展览 and $($allTypesUpdates.Count) cells on 全部类型"
